# Add Mastercard, Amex and Discover test card numbers / exp dates / cvv
# for the Drybar US gold test data sheet (rows 9-11, columns U/V/W).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Row 9 - CCMastercard
$ws.Range("U9").Value = "'5555555555554444"
$ws.Range("V9").Value = "'06/28"
$ws.Range("W9").Value = 123

# Row 10 - CCAmexcard
$ws.Range("U10").Value = "'378282246310005"
$ws.Range("V10").Value = "'06/28"
$ws.Range("W10").Value = 123

# Row 11 - CCDiscovercard
$ws.Range("U11").Value = "'6011111111111117"
$ws.Range("V11").Value = "'06/29"
$ws.Range("W11").Value = 123

# Reflect updated selection state, matching the saved view.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("X10").Select()
